$wb = $excel.ActiveWorkbook

# 1. Update batch size on Charge 08 sheet
$wsCharge = $wb.Worksheets.Item("Charge 08")
$wsCharge.Range("F9").Value2 = 106.5

# 2. Sort stock sheet "A" by column A descending (confirms which stock was used)
$wsA = $wb.Worksheets.Item("A")
$rangeToSort = $wsA.Range("A2:C51")
$keyRange = $wsA.Range("A1")
$wsA.Sort.SortFields.Clear()
$wsA.Sort.SortFields.Add($keyRange, 0, 2, 0, 0) | Out-Null
$wsA.Sort.SetRange($rangeToSort)
$wsA.Sort.Header = 0
$wsA.Sort.Apply()

# The two smallest (last two after descending sort) are unconfirmed / unused stock
$wsA.Cells.Item(8, 2).Value2 = 0
$wsA.Cells.Item(9, 2).Value2 = 0

# 3. Update stock sheet "B" values/confirmations
$wsB = $wb.Worksheets.Item("B")
$wsB.Cells.Item(2, 1).Value2 = 12.114
$wsB.Cells.Item(5, 1).Value2 = 12.664
$wsB.Cells.Item(5, 2).Value2 = 0
$wsB.Cells.Item(6, 1).Value2 = 15.063
$wsB.Cells.Item(6, 2).Value2 = 0
$wsB.Cells.Item(7, 2).Value2 = 0

Write-Host "Charge08 F9:" $wsCharge.Range("F9").Value2
Write-Host "Charge08 G9:" $wsCharge.Range("G9").Value2
Write-Host "Charge08 K9:" $wsCharge.Range("K9").Value2
Write-Host "Charge08 N9:" $wsCharge.Range("N9").Value2
Write-Host "Charge08 F4:" $wsCharge.Range("F4").Value2
Write-Host "Charge08 B21:" $wsCharge.Range("B21").Value2
Write-Host "Charge08 E21:" $wsCharge.Range("E21").Value2

$wsLab = $wb.Worksheets.Item("Lab Print")
Write-Host "LabPrint A3 v:" $wsLab.Range("A3").Value2
Write-Host "LabPrint B6:" $wsLab.Range("B6").Value2
Write-Host "LabPrint C6:" $wsLab.Range("C6").Value2
Write-Host "LabPrint G6:" $wsLab.Range("G6").Value2
Write-Host "LabPrint G21:" $wsLab.Range("G21").Value2
